$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet used to have a spare, blank, styled row at the very top
# (row 1) above the real header row. Deleting it shifts the header
# (old row 2: Name/Math/English/Science/History/Computer) up to row 1
# and every student row up by one, leaving the last former row (7)
# gone and the used range at A1:F6 - matching the target layout.
$ws.Rows.Item(1).Delete()

# Reselect so the active cell matches the new layout (first data row).
$ws.Range("A2").Select()
